# Add-files-via-upload edit: append the "March 30, 2020" snapshot to every
# tracker sheet (new daily case counts), mirroring the previous day's block.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "NYStateCaseTracker2": 56 new county rows (713-768), copied format
# from the prior day's identical 56-row block (657-712) so the alternating
# zebra-stripe styles (and the "big number" comma format for counts >= 1000)
# line up automatically, then the County / ActiveCases / DateTime values are
# written on top.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("NYStateCaseTracker2")

$src1 = $ws1.Range("A657:C712")
$dst1 = $ws1.Range("A713:C768")
$src1.Copy()
$dst1.PasteSpecial(-4122)  # xlPasteFormats

$countyRows = @(
  @('Albany',217),
  @('Allegany',7),
  @('Broome',35),
  @('Cattaraugus',6),
  @('Cayuga',3),
  @('Chautauqua',5),
  @('Chemung',15),
  @('Chenango',17),
  @('Clinton',17),
  @('Columbia',26),
  @('Cortland',8),
  @('Delaware',11),
  @('Dutchess',392),
  @('Erie',376),
  @('Essex',4),
  @('Franklin',6),
  @('Fulton',1),
  @('Genesee',9),
  @('Greene',10),
  @('Hamilton',2),
  @('Herkimer',12),
  @('Jefferson',11),
  @('Lewis',2),
  @('Livingston',12),
  @('Madison',34),
  @('Monroe',242),
  @('Montgomery',6),
  @('Nassau',7344),
  @('Niagara',41),
  @('New York City',37453),
  @('Oneida',34),
  @('Onondaga',180),
  @('Ontario',20),
  @('Orange',1435),
  @('Orleans',4),
  @('Oswego',14),
  @('Otsego',17),
  @('Putnam',167),
  @('Rensselaer',40),
  @('Rockland',2511),
  @('Saratoga',105),
  @('Schenectady',80),
  @('Schoharie',6),
  @('Schuyler',2),
  @('St. Lawrence',13),
  @('Steuben',19),
  @('Suffolk',5791),
  @('Sullivan',101),
  @('Tioga',4),
  @('Tompkins',66),
  @('Ulster',190),
  @('Warren',18),
  @('Washington',7),
  @('Wayne',15),
  @('Westchester',9326),
  @('Wyoming',8)
)

$r = 713
foreach ($row in $countyRows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = "March 30, 2020 4:00PM"
    $r = $r + 1
}

[void]$ws1.Range("D768").Select()

# ---------------------------------------------------------------------------
# Sheet "NYCBoroughs": new 5:00PM snapshot, rows 44-49
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("NYCBoroughs")

$boroughRows = @(
  @('Bronx',6925,1880,215),
  @('Brooklyn',10171,1661,216),
  @('Manhattan',6060,1075,119),
  @('Queens',12756,2650,305),
  @('Staten Island',2140,465,58),
  @('Unknown',35,10,1)
)

$r = 44
foreach ($row in $boroughRows) {
    $ws4.Cells.Item($r, 1).Value = $row[0]
    $ws4.Cells.Item($r, 2).Value = $row[1]
    $ws4.Cells.Item($r, 3).Value = $row[2]
    $ws4.Cells.Item($r, 4).Value = $row[3]
    $ws4.Cells.Item($r, 5).Value = "March 30, 2020 5:00PM"
    $r = $r + 1
}

[void]$ws4.Range("F55").Select()

# ---------------------------------------------------------------------------
# Sheet "NYCCaseDemographics": fill in the missing C43 (Hospitalized=0) on
# the previous snapshot's "Unknown" row, then add the new 5:00PM snapshot,
# rows 44-49 (row 49 "Unknown" has no Deaths/D value, same as prior weeks).
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("NYCCaseDemographics")

$ws5.Cells.Item(43, 3).Value = 0

$ageRows = @(
  @('0 to 17',714,72,1),
  @('18 to 44',16028,1448,54),
  @('45 to 64',13344,2887,216),
  @('65 to 74',4496,1612,215),
  @('75 >',3410,1722,428)
)

$r = 44
foreach ($row in $ageRows) {
    $ws5.Cells.Item($r, 1).Value = $row[0]
    $ws5.Cells.Item($r, 2).Value = $row[1]
    $ws5.Cells.Item($r, 3).Value = $row[2]
    $ws5.Cells.Item($r, 4).Value = $row[3]
    $ws5.Cells.Item($r, 5).Value = "March 30, 2020 5:00PM"
    $r = $r + 1
}

$ws5.Cells.Item(49, 1).Value = "Unknown"
$ws5.Cells.Item(49, 2).Value = 95
$ws5.Cells.Item(49, 3).Value = 0
$ws5.Cells.Item(49, 5).Value = "March 30, 2020 5:00PM"

[void]$ws5.Range("C50").Select()

# ---------------------------------------------------------------------------
# Sheet "NYCCasesSex": new 5:00PM snapshot, rows 26-28
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("NYCCasesSex")

$sexRows = @(
  @('Female',16920,3130,334),
  @('Male',21120,4610,569),
  @('Unknown',47,1,1)
)

$r = 26
foreach ($row in $sexRows) {
    $ws6.Cells.Item($r, 1).Value = $row[0]
    $ws6.Cells.Item($r, 2).Value = $row[1]
    $ws6.Cells.Item($r, 3).Value = $row[2]
    $ws6.Cells.Item($r, 4).Value = $row[3]
    $ws6.Cells.Item($r, 5).Value = "March 30, 2020 5:00PM"
    $r = $r + 1
}

# Re-select sheet6 last so it remains the visible/active tab, matching the
# workbook's activeTab, with the updated selection.
[void]$ws6.Range("D29").Select()
